$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) as text before writing, so numeric-looking
# strings like "52.058.08" / "0.999" / "52.70" keep their exact literal
# formatting instead of being auto-coerced to numbers by Excel.
$priceRange = $ws.Range('D2:D51')
$priceRange.NumberFormat = '@'

$ws.Range('D2').Value = '52.058.08'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.867.07'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '350.81'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '112.09'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('E7').Value = '  +1.37%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = '40.18'
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').Value = '0.135'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '0.0852'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').Value = '20.06'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '7.81'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '3.318.76'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('E16').Value = '  +6.82%  '
$ws.Range('D17').Value = '2.856.25'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').Value = '52.046.79'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = '3.35'
$ws.Range('E19').Value = '  +7.94%  '
$ws.Range('D20').Value = '7.65'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').Value = '13.60'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').Value = '70.77'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').Value = '269.16'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').Value = '2.78'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = '26.30'
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '0.166'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '10.62'
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '38.87'
$ws.Range('E30').Value = '  +4.08%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '2.26'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '6.25'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '5.97'
$ws.Range('E33').Value = '  +6.75%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '52.70'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0920'
$ws.Range('E35').Value = '  +9.30%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '0.0458'
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  +5.99%  '
$ws.Range('D39').Value = '18.54'
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  +2.68%  '
$ws.Range('D41').Value = '2.59'
$ws.Range('E41').Value = '  +2.98%  '
$ws.Range('D42').Value = '0.117'
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = '121.41'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('D44').Value = '22.29'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '2.20'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '3.57'
$ws.Range('E46').Value = '  +7.70%  '
$ws.Range('D47').Value = '2.178.16'
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('D48').Value = '2.49'
$ws.Range('E48').Value = '  +7.05%  '
$ws.Range('D49').Value = '0.247'
$ws.Range('E49').Value = '  +10.74%  '
$ws.Range('D50').Value = '0.953'
$ws.Range('E50').Value = '  +4.65%  '
$ws.Range('D51').Value = '5.51'
$ws.Range('E51').Value = '  +0.67%  '

# Restore the default cell style on the Price column now that the text
# values are set, so no stray number-format style lingers on the cells.
$priceRange.Style = 'Normal'

Write-Host "Applied cryptos update"
